$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows: c_office_id, c_office_chn, c_dy, c_office_trans, c_office_pinyin, c_source
$data = @(
    @(803803, "布政使", "20", "Provincial Administration Commissioner", "bu zheng shi", "67822"),
    @(803804, "書吏", "20", "Clerk", "shu li", "67822"),
    @(803805, "經承", "20", "Assignee", "jing cheng", "67822"),
    @(803806, "教坊司左司樂", "20", "Left Music Director", "jiao fang si zuo si le", "67822"),
    @(803807, "教坊司右司樂", "20", "Right Music Director", "jiao fang si you si le", "67822"),
    @(803808, "教坊司奉鑾", "20", "provider of carriage bells", "jiao fang si feng luan", "67822"),
    @(803809, "翊國公", "19", "Duke of Yi", "yi guo gong", "67822"),
    @(803810, "通議大夫", "20", "Grand Master for Thorough Council", "tong yi dai fu", "67822")
)

$rowIndex = 2
foreach ($entry in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $entry[0]
    $ws.Cells.Item($rowIndex, 2).Value = $entry[1]
    $ws.Cells.Item($rowIndex, 3).Value = $entry[2]
    $ws.Cells.Item($rowIndex, 4).Value = $entry[3]
    $ws.Cells.Item($rowIndex, 5).Value = $entry[4]
    $ws.Cells.Item($rowIndex, 6).Value = $entry[5]
    $rowIndex++
}
